$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 191
$ws1.Range("F4").Value = 385
$ws1.Range("F5").Value = 430
$ws1.Range("F7").Value = 2437
$ws1.Range("F8").Value = 421
$ws1.Range("F9").Value = 6386
$ws1.Range("F10").Value = 170
$ws1.Range("F11").Value = 413

$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F4").Value = 1

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 191
$ws4.Range("F4").Value = 385
$ws4.Range("F5").Value = 430
$ws4.Range("F9").Value = 2437
$ws4.Range("F10").Value = 421
$ws4.Range("F11").Value = 6386
$ws4.Range("F12").Value = 170
$ws4.Range("F13").Value = 413
$ws4.Range("F14").Value = 1
